$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 23 de Abril de 2020 a las 01:22"

# 2) Update Estados Unidos row (row 4)
$ws.Range("B4").Value = 846692
$ws.Range("C4").Value = 27948
$ws.Range("D4").Value = 83921
$ws.Range("E4").Value = 715234
$ws.Range("F4").Value = 14014
$ws.Range("G4").Value = 2219
$ws.Range("H4").Value = 47537

# 3) Update Canada row (row 16)
$ws.Range("B16").Value = 40190
$ws.Range("C16").Value = 1768
$ws.Range("D16").Value = 13986
$ws.Range("E16").Value = 24230
$ws.Range("F16").Value = 557
$ws.Range("G16").Value = 140
$ws.Range("H16").Value = 1974

# 4) Noruega moves above Bielorrusia in ranking (row 39 becomes Noruega with
#    updated figures, row 40 becomes Bielorrusia carrying the old Bielorrusia
#    figures forward unchanged). The country-name cells (column A) keep
#    pointing at the same shared-string slots; only the shared string text
#    assigned to those slots is swapped, which we do by re-assigning the
#    country name values directly.
$ws.Range("A39").Value = "Noruega"
$ws.Range("B39").Value = 7338
$ws.Range("C39").Value = 97
$ws.Range("D39").Value = 32
$ws.Range("E39").Value = 7119
$ws.Range("F39").Value = 54
$ws.Range("G39").Value = 5
$ws.Range("H39").Value = 187

$ws.Range("A40").Value = "Bielorrusia"
$ws.Range("B40").Value = 7281
$ws.Range("C40").Value = 558
$ws.Range("D40").Value = 769
$ws.Range("E40").Value = 6454
$ws.Range("F40").Value = 92
$ws.Range("G40").Value = 3
$ws.Range("H40").Value = 58
